$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidentiality/date notice shared string (cell A40)
$ws.Range("A40").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-37
$ws.Range("D2").Value = 0.03083093971480531
$ws.Range("E2").Value = -0.004664574675600108
$ws.Range("D3").Value = 0.03389939088564103
$ws.Range("E3").Value = 0.006313945224045936
$ws.Range("D4").Value = 0.03371579556633179
$ws.Range("E4").Value = -0.01340571681808111
$ws.Range("D5").Value = 0.06868220584079653
$ws.Range("E5").Value = 0.001174755803076133
$ws.Range("D6").Value = 0.03010888528504994
$ws.Range("E6").Value = 0.006699419383653593
$ws.Range("D7").Value = 0.01582094843138064
$ws.Range("E7").Value = 0.01820371157387735
$ws.Range("D8").Value = 0.03200890073096338
$ws.Range("E8").Value = -0.0005309806804723038
$ws.Range("D9").Value = 0.03169699413559569
$ws.Range("E9").Value = -0.0118200941600719
$ws.Range("D10").Value = 0.05049151436628704
$ws.Range("E10").Value = 0.01344603092402163
$ws.Range("D11").Value = 0.02739118895619057
$ws.Range("E11").Value = 0.009498353300558593
$ws.Range("D12").Value = 0.01461926757234132
$ws.Range("E12").Value = -0.04042210695761039
$ws.Range("D13").Value = 0.01592964881371836
$ws.Range("E13").Value = 0.001758705592683762
$ws.Range("D14").Value = 0.01475785122130109
$ws.Range("E14").Value = -0.01683203401842659
$ws.Range("D15").Value = 0.007018458707022027
$ws.Range("E15").Value = -0.01860130927670445
$ws.Range("D16").Value = 0.007267610442483397
$ws.Range("E16").Value = -0.0234375
$ws.Range("D17").Value = 0.03190300190487148
$ws.Range("E17").Value = -0.006515859355790576
$ws.Range("D18").Value = 0.02800043405444768
$ws.Range("E18").Value = -0.004695869102648631
$ws.Range("D19").Value = 0.03010664404005329
$ws.Range("E19").Value = 0.007556018759770744
$ws.Range("D20").Value = 0.03287999795299623
$ws.Range("E20").Value = -0.01266721576869534
$ws.Range("D21").Value = 0.04893179461903224
$ws.Range("E21").Value = 0.008305692223719463
$ws.Range("D22").Value = 0.02808877646139913
$ws.Range("E22").Value = -0.004202351188893116
$ws.Range("D23").Value = 0.02960778025788139
$ws.Range("E23").Value = 0.006844346317615546
$ws.Range("D24").Value = 0.02805646517936404
$ws.Range("E24").Value = 0.0162962075369959
$ws.Range("D25").Value = 0.01173067631248708
$ws.Range("E25").Value = -0.02942300343905235
$ws.Range("D26").Value = 0.01250054396883773
$ws.Range("E26").Value = -0.007530255490811211
$ws.Range("D27").Value = 0.02882016941197389
$ws.Range("E27").Value = -0.003279155973766712
$ws.Range("D28").Value = 0.02828750018443579
$ws.Range("E28").Value = -0.001162053665751062
$ws.Range("D29").Value = 0.03148725095799216
$ws.Range("E29").Value = 0.00544522741832143
$ws.Range("D30").Value = 0.0331982547425211
$ws.Range("E30").Value = -0.0003094250881862104
$ws.Range("D31").Value = 0.03025680745482911
$ws.Range("E31").Value = -0.02017283950617288
$ws.Range("D32").Value = 0.02856840289068309
$ws.Range("E32").Value = 0.01032949790794957
$ws.Range("D33").Value = 0.03008311096758843
$ws.Range("E33").Value = 0.00395480225988698
$ws.Range("D34").Value = 0.03126835602998562
$ws.Range("E34").Value = -0.009079180006689946
$ws.Range("D35").Value = 0.02907081531076636
$ws.Range("E35").Value = -0.003405075489881115
$ws.Range("D36").Value = 0.03291361662794605
$ws.Range("E36").Value = 0.007467725918570034
$ws.Range("E37").Value = -0.0006951594897955937
